$d = $word.ActiveDocument

$ok = $d.Content.Find.Execute("2025-08-18 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-19 Tuesday", 2)
if (-not $ok) { Write-Host "FAILED to replace: 2025-08-18 Monday" }
$ok = $d.Content.Find.Execute("52×94=4888", $true, $false, $false, $false, $false, $true, 1, $false, "79×82=6478", 2)
if (-not $ok) { Write-Host "FAILED to replace: 52×94=4888" }
$ok = $d.Content.Find.Execute("76×58=4408", $true, $false, $false, $false, $false, $true, 1, $false, "98×63=6174", 2)
if (-not $ok) { Write-Host "FAILED to replace: 76×58=4408" }
$ok = $d.Content.Find.Execute("62×82=5084", $true, $false, $false, $false, $false, $true, 1, $false, "85×58=4930", 2)
if (-not $ok) { Write-Host "FAILED to replace: 62×82=5084" }
$ok = $d.Content.Find.Execute("25×42=1050", $true, $false, $false, $false, $false, $true, 1, $false, "19×27=513", 2)
if (-not $ok) { Write-Host "FAILED to replace: 25×42=1050" }
$ok = $d.Content.Find.Execute("97×47=4559", $true, $false, $false, $false, $false, $true, 1, $false, "35×61=2135", 2)
if (-not $ok) { Write-Host "FAILED to replace: 97×47=4559" }
$ok = $d.Content.Find.Execute("49×57=2793", $true, $false, $false, $false, $false, $true, 1, $false, "15×59=885", 2)
if (-not $ok) { Write-Host "FAILED to replace: 49×57=2793" }
$ok = $d.Content.Find.Execute("34×36=1224", $true, $false, $false, $false, $false, $true, 1, $false, "47×97=4559", 2)
if (-not $ok) { Write-Host "FAILED to replace: 34×36=1224" }
$ok = $d.Content.Find.Execute("20×30=600", $true, $false, $false, $false, $false, $true, 1, $false, "51×49=2499", 2)
if (-not $ok) { Write-Host "FAILED to replace: 20×30=600" }
$ok = $d.Content.Find.Execute("29×46=1334", $true, $false, $false, $false, $false, $true, 1, $false, "27×30=810", 2)
if (-not $ok) { Write-Host "FAILED to replace: 29×46=1334" }
$ok = $d.Content.Find.Execute("38×14=532", $true, $false, $false, $false, $false, $true, 1, $false, "76×84=6384", 2)
if (-not $ok) { Write-Host "FAILED to replace: 38×14=532" }
$ok = $d.Content.Find.Execute("88×61=5368", $true, $false, $false, $false, $false, $true, 1, $false, "60×84=5040", 2)
if (-not $ok) { Write-Host "FAILED to replace: 88×61=5368" }
$ok = $d.Content.Find.Execute("65×89=5785", $true, $false, $false, $false, $false, $true, 1, $false, "11×82=902", 2)
if (-not $ok) { Write-Host "FAILED to replace: 65×89=5785" }
$ok = $d.Content.Find.Execute("74×68=5032", $true, $false, $false, $false, $false, $true, 1, $false, "34×86=2924", 2)
if (-not $ok) { Write-Host "FAILED to replace: 74×68=5032" }
$ok = $d.Content.Find.Execute("81×98=7938", $true, $false, $false, $false, $false, $true, 1, $false, "17×51=867", 2)
if (-not $ok) { Write-Host "FAILED to replace: 81×98=7938" }
$ok = $d.Content.Find.Execute("64×96=6144", $true, $false, $false, $false, $false, $true, 1, $false, "50×94=4700", 2)
if (-not $ok) { Write-Host "FAILED to replace: 64×96=6144" }
$ok = $d.Content.Find.Execute("80×71=5680", $true, $false, $false, $false, $false, $true, 1, $false, "32×85=2720", 2)
if (-not $ok) { Write-Host "FAILED to replace: 80×71=5680" }
$ok = $d.Content.Find.Execute("45×79=3555", $true, $false, $false, $false, $false, $true, 1, $false, "63×62=3906", 2)
if (-not $ok) { Write-Host "FAILED to replace: 45×79=3555" }
$ok = $d.Content.Find.Execute("23×93=2139", $true, $false, $false, $false, $false, $true, 1, $false, "55×96=5280", 2)
if (-not $ok) { Write-Host "FAILED to replace: 23×93=2139" }
$ok = $d.Content.Find.Execute("37×92=3404", $true, $false, $false, $false, $false, $true, 1, $false, "20×93=1860", 2)
if (-not $ok) { Write-Host "FAILED to replace: 37×92=3404" }
$ok = $d.Content.Find.Execute("21×33=693", $true, $false, $false, $false, $false, $true, 1, $false, "31×94=2914", 2)
if (-not $ok) { Write-Host "FAILED to replace: 21×33=693" }
$ok = $d.Content.Find.Execute("41×97=3977", $true, $false, $false, $false, $false, $true, 1, $false, "21×47=987", 2)
if (-not $ok) { Write-Host "FAILED to replace: 41×97=3977" }
$ok = $d.Content.Find.Execute("56×79=4424", $true, $false, $false, $false, $false, $true, 1, $false, "22×63=1386", 2)
if (-not $ok) { Write-Host "FAILED to replace: 56×79=4424" }
$ok = $d.Content.Find.Execute("77×65=5005", $true, $false, $false, $false, $false, $true, 1, $false, "60×54=3240", 2)
if (-not $ok) { Write-Host "FAILED to replace: 77×65=5005" }
$ok = $d.Content.Find.Execute("30×38=1140", $true, $false, $false, $false, $false, $true, 1, $false, "72×82=5904", 2)
if (-not $ok) { Write-Host "FAILED to replace: 30×38=1140" }
$ok = $d.Content.Find.Execute("43×88=3784", $true, $false, $false, $false, $false, $true, 1, $false, "20×44=880", 2)
if (-not $ok) { Write-Host "FAILED to replace: 43×88=3784" }
